# New weekly price record: insert a row at 115 (pushing existing rows 115-176
# down to 116-177) and populate it with the new week's data for
# "Terminal Hortofrutícola Agro Chillán" / Zapallo italiano.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 115..176 down to 116..177, leaving row 115 blank for the new record.
$ws.Rows("115:115").Insert()

# Fill in the new row 115 with the new weekly record.
$ws.Range("A115").Value = 7
$ws.Range("B115").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C115").Value = "Ñuble"
$ws.Range("D115").Value = 44572
$ws.Range("E115").Value = 16
$ws.Range("F115").Value = 100112032
$ws.Range("G115").Value = "Zapallo italiano"
$ws.Range("H115").Value = "Sin especificar"
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 120
$ws.Range("K115").Value = 6500
$ws.Range("L115").Value = 7000
$ws.Range("M115").Value = 6750
$ws.Range("N115").Value = "$/caja 60 unidades"
$ws.Range("O115").Value = "Región del Maule"
$ws.Range("P115").Value = 112
$ws.Range("Q115").Value = 60
$ws.Range("R115").Value = "Hortaliza"
